# Daily attendance processing - 2026-01-06 23:02:34
#
# In the "Recorded By" column (G), when the recorded-by value lists
# "System" as the second entry after the user's e-mail address
# (e.g. "dnasr281@gmail.com, System"), swap the order so "System" is
# listed first (e.g. "System, dnasr281@gmail.com").
#
# Uses Find/FindNext (like Excel's own "Find All") so only the cells that
# already contain the text are touched - blank cells in the sheet are left
# completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$searchText = ", System"
$matches = @()

$first = $ws.Cells.Find($searchText)
if ($first) {
    $current = $first
    $continue = $true
    while ($continue) {
        $matches += $current
        $current = $ws.Cells.FindNext($current)
        if (-not $current -or $current.Address() -eq $first.Address()) {
            $continue = $false
        }
    }
}

foreach ($cell in $matches) {
    $val = $cell.Value()
    if ($val -and $val -like "*, System") {
        $emailPart = $val.Substring(0, $val.Length - ", System".Length)
        $cell.Value = "System, " + $emailPart
    }
}
